$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variants = @("Variant","4F","4H","4L","4N","4Q","4V","5M","8C","8Q","12C","12M","71C")
$fitness  = @("fitness",1.0455,0.989,0.9845,1.0115,0.994,0.99775,1.038,0.8985,0.913,0.7385,1.012,1)

for ($i = 0; $i -lt $variants.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $variants[$i]
    $ws.Cells.Item($row, 2).Value = $fitness[$i]
}
